# Added local currency to USD conversion for comparison purposes
#
# The "Data" sheet header row gets two changes:
#   1. Columns AB1/AC1 ("MSRP"/"Reseller Cost") swap places, so the header
#      order becomes ... Cost, Reseller Cost, MSRP, Seats, ...
#   2. Three new USD-equivalent columns are appended after the existing
#      "Seats" column (AD): USD Cost, USD Reseller Cost, USD MSRP.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap MSRP and Reseller Cost header labels (columns AB and AC).
$ws.Range("AB1").Value2 = "Reseller Cost"
$ws.Range("AC1").Value2 = "MSRP"

# Append the new USD comparison columns after "Seats" (column AD).
# Order of assignment controls the order new strings are added to the
# shared string table, so write "USD Cost" then "USD MSRP" then
# "USD Reseller Cost" to match the expected workbook layout, even though
# "USD Reseller Cost" ends up visually in column AF and "USD MSRP" in AG.
$ws.Range("AE1").Value2 = "USD Cost"
$ws.Range("AG1").Value2 = "USD MSRP"
$ws.Range("AF1").Value2 = "USD Reseller Cost"

# Match the saved selection state (active cell moves to AF2).
[void]$ws.Range("AF2").Select()
